$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.1958333333333333
$ws.Range("C2").Value = 0.5583333333333333
$ws.Range("J2").Value = 0.008333333333333333
$ws.Range("P2").Value = 0.1458333333333333
$ws.Range("S2").Value = 0.09166666666666666
$ws.Range("C3").Value = 0.01438848920863309
$ws.Range("J3").Value = 0.05035971223021583
$ws.Range("P3").Value = 0.7338129496402878
$ws.Range("S3").Value = 0.2014388489208633
$ws.Range("J4").Value = 0.07547169811320754
$ws.Range("P4").Value = 0.6415094339622641
$ws.Range("S4").Value = 0.2830188679245283
$ws.Range("S5").Value = 1
$ws.Range("B6").Value = 0.05909090909090909
$ws.Range("D6").Value = 0.01818181818181818
$ws.Range("F6").Value = 0.05454545454545454
$ws.Range("J6").Value = 0.2909090909090909
$ws.Range("O6").Value = 0.02272727272727273
$ws.Range("Q6").Value = 0.1363636363636364
$ws.Range("R6").Value = 0.05
$ws.Range("S6").Value = 0.3681818181818182
$ws.Range("B7").Value = 0.06185567010309279
$ws.Range("D7").Value = 0.04639175257731959
$ws.Range("E7").Value = 0.01030927835051546
$ws.Range("F7").Value = 0.07731958762886598
$ws.Range("J7").Value = 0.1494845360824742
$ws.Range("Q7").Value = 0.1958762886597938
$ws.Range("R7").Value = 0.06701030927835051
$ws.Range("S7").Value = 0.3917525773195876
$ws.Range("B8").Value = 0.07142857142857142
$ws.Range("D8").Value = 0.02678571428571428
$ws.Range("F8").Value = 0.05133928571428571
$ws.Range("J8").Value = 0.1004464285714286
$ws.Range("O8").Value = 0.01339285714285714
$ws.Range("Q8").Value = 0.1808035714285714
$ws.Range("R8").Value = 0.09598214285714286
$ws.Range("S8").Value = 0.4598214285714285
$ws.Range("B9").Value = 0.05947955390334572
$ws.Range("D9").Value = 0.01115241635687732
$ws.Range("F9").Value = 0.07806691449814127
$ws.Range("J9").Value = 0.1263940520446097
$ws.Range("O9").Value = 0.003717472118959108
$ws.Range("Q9").Value = 0.1933085501858736
$ws.Range("R9").Value = 0.07806691449814127
$ws.Range("S9").Value = 0.449814126394052
$ws.Range("B10").Value = 0.08751902587519025
$ws.Range("D10").Value = 0.0213089802130898
$ws.Range("F10").Value = 0.0639269406392694
$ws.Range("J10").Value = 0.1118721461187215
$ws.Range("O10").Value = 0.01598173515981735
$ws.Range("Q10").Value = 0.2267884322678843
$ws.Range("R10").Value = 0.0837138508371385
$ws.Range("S10").Value = 0.3888888888888889
$ws.Range("G11").Value = 0.1547169811320755
$ws.Range("J11").Value = 0.07924528301886792
$ws.Range("K11").Value = 0.1849056603773585
$ws.Range("L11").Value = 0.5773584905660377
$ws.Range("S11").Value = 0.003773584905660377
$ws.Range("G12").Value = 0.7530864197530864
$ws.Range("J12").Value = 0.191358024691358
$ws.Range("K12").Value = 0.006172839506172839
$ws.Range("L12").Value = 0.01851851851851852
$ws.Range("S12").Value = 0.0308641975308642
$ws.Range("G13").Value = 0.8297872340425532
$ws.Range("J13").Value = 0.148936170212766
$ws.Range("S13").Value = 0.02127659574468085
$ws.Range("F15").Value = 0.01809954751131222
$ws.Range("H15").Value = 0.1809954751131222
$ws.Range("I15").Value = 0.1040723981900453
$ws.Range("J15").Value = 0.3574660633484163
$ws.Range("K15").Value = 0.05429864253393665
$ws.Range("M15").Value = 0.004524886877828055
$ws.Range("N15").Value = 0.004524886877828055
$ws.Range("O15").Value = 0.06787330316742081
$ws.Range("S15").Value = 0.2081447963800905
$ws.Range("F16").Value = 0.01807228915662651
$ws.Range("H16").Value = 0.1204819277108434
$ws.Range("I16").Value = 0.08433734939759036
$ws.Range("J16").Value = 0.4578313253012048
$ws.Range("K16").Value = 0.1385542168674699
$ws.Range("M16").Value = 0.01204819277108434
$ws.Range("N16").Value = 0.006024096385542169
$ws.Range("O16").Value = 0.04819277108433735
$ws.Range("S16").Value = 0.1144578313253012
$ws.Range("F17").Value = 0.0163265306122449
$ws.Range("H17").Value = 0.1714285714285714
$ws.Range("I17").Value = 0.1224489795918367
$ws.Range("J17").Value = 0.4224489795918367
$ws.Range("K17").Value = 0.07346938775510205
$ws.Range("M17").Value = 0.01224489795918367
$ws.Range("N17").Value = 0.002040816326530612
$ws.Range("O17").Value = 0.0653061224489796
$ws.Range("S17").Value = 0.1142857142857143
$ws.Range("F18").Value = 0.01530612244897959
$ws.Range("H18").Value = 0.1377551020408163
$ws.Range("I18").Value = 0.1326530612244898
$ws.Range("J18").Value = 0.4336734693877551
$ws.Range("K18").Value = 0.07653061224489796
$ws.Range("M18").Value = 0.01530612244897959
$ws.Range("O18").Value = 0.0663265306122449
$ws.Range("S18").Value = 0.1224489795918367
$ws.Range("F19").Value = 0.01622418879056047
$ws.Range("H19").Value = 0.2042772861356932
$ws.Range("I19").Value = 0.107669616519174
$ws.Range("J19").Value = 0.3657817109144543
$ws.Range("K19").Value = 0.09144542772861357
$ws.Range("M19").Value = 0.02654867256637168
$ws.Range("O19").Value = 0.06932153392330384
$ws.Range("S19").Value = 0.1187315634218289
